$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 19.04.2025"

# --- Row 6 ---
$ws.Range("B6").Value = "21.04."
$ws.Range("C6").Value = "22.04."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-16594652"
$ws.Range("E6").Value = "54,11-"

# --- Row 7 ---
$ws.Range("B7").Value = "25.04."
$ws.Range("C7").Value = "26.04."
$ws.Range("D7").Value = "BURGER KING Eisleben"
$ws.Range("E7").Value = "17,57-"

# --- Row 8 ---
$ws.Range("B8").Value = "27.04."
$ws.Range("C8").Value = "28.04."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 85356004"
$ws.Range("E8").Value = "85,92-"

# --- Row 9: the 4th transaction is removed -> row becomes empty ---
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 30.04.2025"
$ws.Range("E12").Value = "157,60-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 06.05.2025"
